# Updates cryptos list values (price / volume change / swapped rows) per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "26.555.44"
Set-TextCell 2 5 "  +1.10%  "

Set-TextCell 3 4 "1.813.23"
Set-TextCell 3 5 "  +1.24%  "

Set-TextCell 4 4 "1.007"
Set-TextCell 4 5 "  -0.17%  "

Set-TextCell 5 4 "1.006"
Set-TextCell 5 5 "  -0.16%  "

Set-TextCell 6 4 "305.33"
Set-TextCell 6 5 "  -0.53%  "

Set-TextCell 7 4 "0.4649"
Set-TextCell 7 5 "  +2.23%  "

Set-TextCell 8 4 "0.3572"
Set-TextCell 8 5 "  -1.02%  "

Set-TextCell 9 5 "  +3.11%  "

Set-TextCell 10 4 "0.07101"
Set-TextCell 10 5 "  +0.57%  "

Set-TextCell 11 4 "0.8985"
Set-TextCell 11 5 "  +3.30%  "

Set-TextCell 12 4 "0.07768"
Set-TextCell 12 5 "  +0.04%  "

Set-TextCell 13 4 "19.32"
Set-TextCell 13 5 "  +0.07%  "

Set-TextCell 14 4 "1.842.48"
Set-TextCell 14 5 "  +3.00%  "

Set-TextCell 15 4 "5.236"
Set-TextCell 15 5 "  -0.48%  "

Set-TextCell 16 4 "6.298"
Set-TextCell 16 5 "  -0.19%  "

Set-TextCell 17 4 "87.53"
Set-TextCell 17 5 "  +3.41%  "

Set-TextCell 18 4 "1.007"
Set-TextCell 18 5 "  -0.20%  "

Set-TextCell 19 4 "0.000008535"
Set-TextCell 19 5 "  +0.65%  "

Set-TextCell 20 5 "  -0.26%  "

Set-TextCell 21 4 "26.591.51"
Set-TextCell 21 5 "  +0.97%  "

Set-TextCell 22 4 "14.14"
Set-TextCell 22 5 "  +0.38%  "

Set-TextCell 23 4 "4.974"
Set-TextCell 23 5 "  +0.16%  "

Set-TextCell 24 4 "10.52"
Set-TextCell 24 5 "  +0.48%  "

Set-TextCell 25 4 "1.921"
Set-TextCell 25 5 "  -2.71%  "

Set-TextCell 26 4 "151.89"
Set-TextCell 26 5 "  -0.23%  "

Set-TextCell 27 4 "17.84"
Set-TextCell 27 5 "  +0.40%  "

Set-TextCell 28 4 "1.998"
Set-TextCell 28 5 "  -1.46%  "

Set-TextCell 29 4 "112.83"
Set-TextCell 29 5 "  +0.61%  "

Set-TextCell 30 5 "  -0.51%  "

Set-TextCell 31 4 "0.08724"
Set-TextCell 31 5 "  +0.96%  "

Set-TextCell 32 4 "3.114"
Set-TextCell 32 5 "  +2.95%  "

Set-TextCell 33 2 "RenderToken"
Set-TextCell 33 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 33 4 "2.745"
Set-TextCell 33 5 "  +4.46%  "

Set-TextCell 34 2 "ImmutableX"
Set-TextCell 34 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 34 4 "0.7302"
Set-TextCell 34 5 "  +2.52%  "

Set-TextCell 35 4 "4.419"
Set-TextCell 35 5 "  -0.28%  "

Set-TextCell 36 4 "1.119"
Set-TextCell 36 5 "  +1.91%  "

Set-TextCell 37 4 "1.073"
Set-TextCell 37 5 "  -0.23%  "

Set-TextCell 38 4 "0.01926"
Set-TextCell 38 5 "  -0.42%  "

Set-TextCell 39 4 "2.913"
Set-TextCell 39 5 "  +1.78%  "

Set-TextCell 40 4 "0.05082"
Set-TextCell 40 5 "  +0.07%  "

Set-TextCell 41 4 "0.5047"
Set-TextCell 41 5 "  +2.82%  "

Set-TextCell 42 4 "6.799"
Set-TextCell 42 5 "  -1.08%  "

Set-TextCell 43 4 "0.1492"
Set-TextCell 43 5 "  -1.71%  "

Set-TextCell 44 4 "7.954"
Set-TextCell 44 5 "  +0.05%  "

Set-TextCell 45 4 "0.4674"
Set-TextCell 45 5 "  +2.51%  "

Set-TextCell 46 5 "  -0.20%  "

Set-TextCell 47 4 "9.988"
Set-TextCell 47 5 "  +1.11%  "

Set-TextCell 48 4 "98.03"
Set-TextCell 48 5 "  -1.67%  "

Set-TextCell 49 4 "1.564"
Set-TextCell 49 5 "  -0.81%  "

Set-TextCell 50 4 "0.06009"
Set-TextCell 50 5 "  +1.12%  "

Set-TextCell 51 4 "63.54"
Set-TextCell 51 5 "  +0.26%  "
